# Updated cryptos list data (price + 1h volume change) to match the latest
# coinranking.com snapshot. Also fixes the WEMIXToken/RenderToken row order
# swap (rows 46-47), which flipped position in this update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values below are written as plain text (matching the workbook's
# original inline-string cell type). Some "Price" column values look like
# plain decimal numbers (e.g. "323.38"), and Excel's COM automation will
# silently coerce those into numeric cells unless the cell is pre-formatted
# as Text ("@"). We apply that number format only to the handful of cells
# that need it, then immediately clear the format again afterwards so the
# cell keeps its original (unstyled) appearance while the value stays text.
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).ClearFormats()
}

# Row 2
$ws.Range("D2").Value = '28.090.72'
$ws.Range("E2").Value = '  -1.99%  '

# Row 3
$ws.Range("D3").Value = '1.836.49'
$ws.Range("E3").Value = '  -0.83%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
Set-TextValue $ws "D5" '323.38'
$ws.Range("E5").Value = '  -3.36%  '

# Row 6
$ws.Range("E6").Value = '  -0.05%  '

# Row 7
Set-TextValue $ws "D7" '0.4629'
$ws.Range("E7").Value = '  -0.48%  '

# Row 8
Set-TextValue $ws "D8" '0.3879'
$ws.Range("E8").Value = '  -0.76%  '

# Row 9
$ws.Range("E9").Value = '  -0.64%  '

# Row 10
Set-TextValue $ws "D10" '0.9633'
$ws.Range("E10").Value = '  -2.16%  '

# Row 11
Set-TextValue $ws "D11" '21.97'
$ws.Range("E11").Value = '  -1.52%  '

# Row 12
$ws.Range("D12").Value = '1.806.23'
$ws.Range("E12").Value = '  -2.64%  '

# Row 13
Set-TextValue $ws "D13" '5.698'
$ws.Range("E13").Value = '  -2.62%  '

# Row 14
Set-TextValue $ws "D14" '6.931'
$ws.Range("E14").Value = '  -1.12%  '

# Row 15
Set-TextValue $ws "D15" '0.06829'
$ws.Range("E15").Value = '  -0.39%  '

# Row 16
Set-TextValue $ws "D16" '88.54'
$ws.Range("E16").Value = '  +0.89%  '

# Row 18
Set-TextValue $ws "D18" '0.000009960'
$ws.Range("E18").Value = '  -1.26%  '

# Row 19
$ws.Range("E19").Value = '  -2.56%  '

# Row 20
Set-TextValue $ws "D20" '1.001'
$ws.Range("E20").Value = '  -0.12%  '

# Row 21
$ws.Range("D21").Value = '28.112.29'
$ws.Range("E21").Value = '  -1.96%  '

# Row 22
Set-TextValue $ws "D22" '5.326'
$ws.Range("E22").Value = '  -1.39%  '

# Row 23
Set-TextValue $ws "D23" '11.01'
$ws.Range("E23").Value = '  -2.68%  '

# Row 24
Set-TextValue $ws "D24" '2.098'
$ws.Range("E24").Value = '  -1.67%  '

# Row 25
$ws.Range("D25").Value = '2.075.16'
$ws.Range("E25").Value = '  -0.63%  '

# Row 26
Set-TextValue $ws "D26" '154.57'
$ws.Range("E26").Value = '  +0.92%  '

# Row 27
$ws.Range("E27").Value = '  -1.56%  '

# Row 28
Set-TextValue $ws "D28" '5.694'
$ws.Range("E28").Value = '  -5.90%  '

# Row 29
Set-TextValue $ws "D29" '1.968'
$ws.Range("E29").Value = '  -3.12%  '

# Row 30
Set-TextValue $ws "D30" '118.05'
$ws.Range("E30").Value = '  +0.51%  '

# Row 31
Set-TextValue $ws "D31" '0.9369'
$ws.Range("E31").Value = '  -4.32%  '

# Row 32
Set-TextValue $ws "D32" '0.09248'
$ws.Range("E32").Value = '  -1.84%  '

# Row 33
Set-TextValue $ws "D33" '5.279'
$ws.Range("E33").Value = '  -1.87%  '

# Row 34
Set-TextValue $ws "D34" '1.324'

# Row 35
Set-TextValue $ws "D35" '3.306'
$ws.Range("E35").Value = '  -5.08%  '

# Row 36
Set-TextValue $ws "D36" '0.05881'
$ws.Range("E36").Value = '  -4.50%  '

# Row 37
Set-TextValue $ws "D37" '0.02129'
$ws.Range("E37").Value = '  -3.27%  '

# Row 38
Set-TextValue $ws "D38" '1.147'
$ws.Range("E38").Value = '  -1.96%  '

# Row 39
Set-TextValue $ws "D39" '7.783'
$ws.Range("E39").Value = '  +2.22%  '

# Row 40
Set-TextValue $ws "D40" '0.5605'
$ws.Range("E40").Value = '  -2.10%  '

# Row 41
Set-TextValue $ws "D41" '9.916'
$ws.Range("E41").Value = '  -2.88%  '

# Row 42
Set-TextValue $ws "D42" '0.1765'
$ws.Range("E42").Value = '  -2.08%  '

# Row 43
Set-TextValue $ws "D43" '0.07260'
$ws.Range("E43").Value = '  +1.72%  '

# Row 44
Set-TextValue $ws "D44" '11.70'
$ws.Range("E44").Value = '  -1.02%  '

# Row 45
Set-TextValue $ws "D45" '0.5286'
$ws.Range("E45").Value = '  -1.99%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws "D46" '1.148'
$ws.Range("E46").Value = '  -7.76%  '

# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws "D47" '2.126'
$ws.Range("E47").Value = '  -10.36%  '

# Row 48
Set-TextValue $ws "D48" '1.831'
$ws.Range("E48").Value = '  -4.08%  '

# Row 49
Set-TextValue $ws "D49" '112.53'
$ws.Range("E49").Value = '  -2.54%  '

# Row 50
Set-TextValue $ws "D50" '1.029'
$ws.Range("E50").Value = '  +0.71%  '

# Row 51
Set-TextValue $ws "D51" '1.001'
$ws.Range("E51").Value = '  -0.07%  '
